$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 2).Value = 40.1
$ws.Cells.Item($row, 3).Value = 57.3
$ws.Cells.Item($row, 4).Value = 60.8
$ws.Cells.Item($row, 5).Value = 65.40000000000001
$ws.Cells.Item($row, 6).Value = 81.40000000000001
$ws.Cells.Item($row, 7).Value = 55.6
$ws.Cells.Item($row, 8).Value = 52.5
$ws.Cells.Item($row, 9).Value = 51.6
$ws.Cells.Item($row, 10).Value = 62
$ws.Cells.Item($row, 11).Value = 61.9
$ws.Cells.Item($row, 12).Value = 56.2
$ws.Cells.Item($row, 13).Value = 3.9
$ws.Cells.Item($row, 14).Value = 60.9
